{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\ntable.values = [\n  [\n    \"7+59=66\",\n    \"58+4=62\",\n    \"37+4=41\",\n    \"96-39=57\",\n    \"57-29=28\"\n  ],\n  [\n    \"8+58=66\",\n    \"72-59=13\",\n    \"83-34=49\",\n    \"19+14=33\",\n    \"76-58=18\"\n  ],\n  [\n    \"54+28=82\",\n    \"5+66=71\",\n    \"28+23=51\",\n    \"18+4=22\",\n    \"28+29=57\"\n  ],\n  [\n    \"18+43=61\",\n    \"38+48=86\",\n    \"55-38=17\",\n    \"8+77=85\",\n    \"67+18=85\"\n  ],\n  [\n    \"51-39=12\",\n    \"35+27=62\",\n    \"89+6=95\",\n    \"32-17=15\",\n    \"93-69=24\"\n  ],\n  [\n    \"33+58=91\",\n    \"69+8=77\",\n    \"4+88=92\",\n    \"85-66=19\",\n    \"63-59=4\"\n  ],\n  [\n    \"49+17=66\",\n    \"56+15=71\",\n    \"95-56=39\",\n    \"12+29=41\",\n    \"94-6=88\"\n  ],\n  [\n    \"79+2=81\",\n    \"18+46=64\",\n    \"83-49=34\",\n    \"81-55=26\",\n    \"28+4=32\"\n  ],\n  [\n    \"61-49=12\",\n    \"43+29=72\",\n    \"65-28=37\",\n    \"15+18=33\",\n    \"57-19=38\"\n  ],\n  [\n    \"23+9=32\",\n    \"51-8=43\",\n    \"96-87=9\",\n    \"67+18=85\",\n    \"29+36=65\"\n  ],\n  [\n    \"26+47=73\",\n    \"83-26=57\",\n    \"64-48=16\",\n    \"29+65=94\",\n    \"75-68=7\"\n  ],\n  [\n    \"70-22=48\",\n    \"28+4=32\",\n    \"50-32=18\",\n    \"38+47=85\",\n    \"37+9=46\"\n  ],\n  [\n    \"4+7=11\",\n    \"57-8=49\",\n    \"22+49=71\",\n    \"94-17=77\",\n    \"39+3=42\"\n  ],\n  [\n    \"83-27=56\",\n    \"65-38=27\",\n    \"71-8=63\",\n    \"53-34=19\",\n    \"37+9=46\"\n  ],\n  [\n    \"88-39=49\",\n    \"29+43=72\",\n    \"24+59=83\",\n    \"53-15=38\",\n    \"7+46=53\"\n  ],\n  [\n    \"89+6=95\",\n    \"26+15=41\",\n    \"25+17=42\",\n    \"39+7=46\",\n    \"90-24=66\"\n  ],\n  [\n    \"39+23=62\",\n    \"21-18=3\",\n    \"51-7=44\",\n    \"62-54=8\",\n    \"73-26=47\"\n  ],\n  [\n    \"13+8=21\",\n    \"63-28=35\",\n    \"36+29=65\",\n    \"80-42=38\",\n    \"19+34=53\"\n  ],\n  [\n    \"40-38=2\",\n    \"15+47=62\",\n    \"72-59=13\",\n    \"7+29=36\",\n    \"9+86=95\"\n  ],\n  [\n    \"66-47=19\",\n    \"63-46=17\",\n    \"52-28=24\",\n    \"13-7=6\",\n    \"68-29=39\"\n  ]\n];\nawait context.sync();\nreturn \"ok\";\n", "ps1": "$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n$tbl.Cell(1,1).Range.Text = \"7+59=66\"\n$tbl.Cell(1,2).Range.Text = \"58+4=62\"\n$tbl.Cell(1,3).Range.Text = \"37+4=41\"\n$tbl.Cell(1,4).Range.Text = \"96-39=57\"\n$tbl.Cell(1,5).Range.Text = \"57-29=28\"\n$tbl.Cell(2,1).Range.Text = \"8+58=66\"\n$tbl.Cell(2,2).Range.Text = \"72-59=13\"\n$tbl.Cell(2,3).Range.Text = \"83-34=49\"\n$tbl.Cell(2,4).Range.Text = \"19+14=33\"\n$tbl.Cell(2,5).Range.Text = \"76-58=18\"\n$tbl.Cell(3,1).Range.Text = \"54+28=82\"\n$tbl.Cell(3,2).Range.Text = \"5+66=71\"\n$tbl.Cell(3,3).Range.Text = \"28+23=51\"\n$tbl.Cell(3,4).Range.Text = \"18+4=22\"\n$tbl.Cell(3,5).Range.Text = \"28+29=57\"\n$tbl.Cell(4,1).Range.Text = \"18+43=61\"\n$tbl.Cell(4,2).Range.Text = \"38+48=86\"\n$tbl.Cell(4,3).Range.Text = \"55-38=17\"\n$tbl.Cell(4,4).Range.Text = \"8+77=85\"\n$tbl.Cell(4,5).Range.Text = \"67+18=85\"\n$tbl.Cell(5,1).Range.Text = \"51-39=12\"\n$tbl.Cell(5,2).Range.Text = \"35+27=62\"\n$tbl.Cell(5,3).Range.Text = \"89+6=95\"\n$tbl.Cell(5,4).Range.Text = \"32-17=15\"\n$tbl.Cell(5,5).Range.Text = \"93-69=24\"\n$tbl.Cell(6,1).Range.Text = \"33+58=91\"\n$tbl.Cell(6,2).Range.Text = \"69+8=77\"\n$tbl.Cell(6,3).Range.Text = \"4+88=92\"\n$tbl.Cell(6,4).Range.Text = \"85-66=19\"\n$tbl.Cell(6,5).Range.Text = \"63-59=4\"\n$tbl.Cell(7,1).Range.Text = \"49+17=66\"\n$tbl.Cell(7,2).Range.Text = \"56+15=71\"\n$tbl.Cell(7,3).Range.Text = \"95-56=39\"\n$tbl.Cell(7,4).Range.Text = \"12+29=41\"\n$tbl.Cell(7,5).Range.Text = \"94-6=88\"\n$tbl.Cell(8,1).Range.Text = \"79+2=81\"\n$tbl.Cell(8,2).Range.Text = \"18+46=64\"\n$tbl.Cell(8,3).Range.Text = \"83-49=34\"\n$tbl.Cell(8,4).Range.Text = \"81-55=26\"\n$tbl.Cell(8,5).Range.Text = \"28+4=32\"\n$tbl.Cell(9,1).Range.Text = \"61-49=12\"\n$tbl.Cell(9,2).Range.Text = \"43+29=72\"\n$tbl.Cell(9,3).Range.Text = \"65-28=37\"\n$tbl.Cell(9,4).Range.Text = \"15+18=33\"\n$tbl.Cell(9,5).Range.Text = \"57-19=38\"\n$tbl.Cell(10,1).Range.Text = \"23+9=32\"\n$tbl.Cell(10,2).Range.Text = \"51-8=43\"\n$tbl.Cell(10,3).Range.Text = \"96-87=9\"\n$tbl.Cell(10,4).Range.Text = \"67+18=85\"\n$tbl.Cell(10,5).Range.Text = \"29+36=65\"\n$tbl.Cell(11,1).Range.Text = \"26+47=73\"\n$tbl.Cell(11,2).Range.Text = \"83-26=57\"\n$tbl.Cell(11,3).Range.Text = \"64-48=16\"\n$tbl.Cell(11,4).Range.Text = \"29+65=94\"\n$tbl.Cell(11,5).Range.Text = \"75-68=7\"\n$tbl.Cell(12,1).Range.Text = \"70-22=48\"\n$tbl.Cell(12,2).Range.Text = \"28+4=32\"\n$tbl.Cell(12,3).Range.Text = \"50-32=18\"\n$tbl.Cell(12,4).Range.Text = \"38+47=85\"\n$tbl.Cell(12,5).Range.Text = \"37+9=46\"\n$tbl.Cell(13,1).Range.Text = \"4+7=11\"\n$tbl.Cell(13,2).Range.Text = \"57-8=49\"\n$tbl.Cell(13,3).Range.Text = \"22+49=71\"\n$tbl.Cell(13,4).Range.Text = \"94-17=77\"\n$tbl.Cell(13,5).Range.Text = \"39+3=42\"\n$tbl.Cell(14,1).Range.Text = \"83-27=56\"\n$tbl.Cell(14,2).Range.Text = \"65-38=27\"\n$tbl.Cell(14,3).Range.Text = \"71-8=63\"\n$tbl.Cell(14,4).Range.Text = \"53-34=19\"\n$tbl.Cell(14,5).Range.Text = \"37+9=46\"\n$tbl.Cell(15,1).Range.Text = \"88-39=49\"\n$tbl.Cell(15,2).Range.Text = \"29+43=72\"\n$tbl.Cell(15,3).Range.Text = \"24+59=83\"\n$tbl.Cell(15,4).Range.Text = \"53-15=38\"\n$tbl.Cell(15,5).Range.Text = \"7+46=53\"\n$tbl.Cell(16,1).Range.Text = \"89+6=95\"\n$tbl.Cell(16,2).Range.Text = \"26+15=41\"\n$tbl.Cell(16,3).Range.Text = \"25+17=42\"\n$tbl.Cell(16,4).Range.Text = \"39+7=46\"\n$tbl.Cell(16,5).Range.Text = \"90-24=66\"\n$tbl.Cell(17,1).Range.Text = \"39+23=62\"\n$tbl.Cell(17,2).Range.Text = \"21-18=3\"\n$tbl.Cell(17,3).Range.Text = \"51-7=44\"\n$tbl.Cell(17,4).Range.Text = \"62-54=8\"\n$tbl.Cell(17,5).Range.Text = \"73-26=47\"\n$tbl.Cell(18,1).Range.Text = \"13+8=21\"\n$tbl.Cell(18,2).Range.Text = \"63-28=35\"\n$tbl.Cell(18,3).Range.Text = \"36+29=65\"\n$tbl.Cell(18,4).Range.Text = \"80-42=38\"\n$tbl.Cell(18,5).Range.Text = \"19+34=53\"\n$tbl.Cell(19,1).Range.Text = \"40-38=2\"\n$tbl.Cell(19,2).Range.Text = \"15+47=62\"\n$tbl.Cell(19,3).Range.Text = \"72-59=13\"\n$tbl.Cell(19,4).Range.Text = \"7+29=36\"\n$tbl.Cell(19,5).Range.Text = \"9+86=95\"\n$tbl.Cell(20,1).Range.Text = \"66-47=19\"\n$tbl.Cell(20,2).Range.Text = \"63-46=17\"\n$tbl.Cell(20,3).Range.Text = \"52-28=24\"\n$tbl.Cell(20,4).Range.Text = \"13-7=6\"\n$tbl.Cell(20,5).Range.Text = \"68-29=39\"\n"}
